$d = $word.ActiveDocument

# Temporarily append an extra trailing paragraph. This prevents paragraph 7
# (Feature Overview body) from being the literal last paragraph of the document
# while we edit it, which avoids an engine quirk that forces xml:space="preserve"
# onto every <w:t> run-child of the last paragraph whenever it is edited.
$tailRng = $d.Content
$tailRng.Collapse(0)
$tailRng.InsertParagraphAfter()

function Replace-Once($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find text not found: $findText"
    }
    $rng.Text = $replaceText
}

# 1. Update the "Generated" timestamp line
Replace-Once "Generated: 2025-09-05 09:00:26" "Generated: 2025-09-05 13:13:47"

# 2. Product Overview - first body paragraph text
$old4 = "This product is a modern, in-house credit card core platform designed specifically for the financial services industry within the United States credit card sector. It offers comprehensive credit card issuance and account management capabilities, supporting both physical and virtual card issuance in real time. The platform delivers end-to-end lifecycle management including credit line administration, product configuration, transaction authorizations, settlements, billing, payments, interest calculation, fee management, rewards, dispute handling, and delinquency workflows. It aims to fully replace legacy third-party systems with a robust, scalable infrastructure enabling full ownership and control over credit issuance and servicing processes."
$new4 = "This product is an in-house credit card core platform designed specifically for the Property & Casualty insurance sector within the financial services industry, targeting the U.S. market. It replaces the existing legacy third-party systems by providing full control over credit issuance, account management, servicing, and related credit lifecycle functions. Delivered as an integrated digital platform, it supports real-time issuance of both virtual and physical credit cards while handling all key processes including credit line management, product setup, authorizations, settlements, billing, payments, interest calculation, fees, rewards, disputes, and delinquency workflows."
Replace-Once $old4 $new4

# 3. Product Overview - second body paragraph text, and remove the trailing
#    "### References" / "No external sources used." block (and its leading break pair)
$old5 = "Strategically, this solution addresses the need for enhanced operational control, agility, and efficiency by transitioning from outdated legacy technologies to a modern core system. It supports sophisticated credit account management including charged-off accounts, enabling accurate lifecycle management aligned with regulatory and business requirements. This transformation allows the enterprise to streamline processes, improve customer experience, reduce dependency risks, and maintain compliance with applicable U.S. financial regulations, positioning the organization for sustainable growth and competitive differentiation in the evolving credit card market." + [char]11 + [char]11 + "### References" + [char]11 + "No external sources used."
$new5 = "Strategically, this platform enables the organization to own the end-to-end credit card ecosystem, enhancing operational efficiency, flexibility, and customer experience while aligning with regulatory and compliance requirements specific to the U.S. insurance and financial sectors. It also introduces modern core infrastructure capable of managing complex account states such as charged-off accounts within trade credit, thereby streamlining transitions triggered by prolonged delinquency, customer death, or bankruptcy. This shift supports the company’s goal of reducing reliance on legacy codebases and outsourced services, positioning it for future scalability and innovation."
Replace-Once $old5 $new5

# 4. Feature Overview - first body paragraph text
$old11 = "The Create a Frontbook Charged-Off Account feature enables the establishment and full lifecycle management of charged-off credit card accounts within a modern in-house credit card core system. This capability is a specialized subset of the broader account creation process, activated when an account transitions from good standing to charged-off status due to triggers such as prolonged delinquency, customer death, or bankruptcy. It supports real-time processing and integration with credit issuance, account management, billing, payments, disputes, and delinquency workflows, ensuring seamless handling of charged-off accounts within the trade credit ecosystem."
$new11 = "The Create a Frontbook Charged-Off Account feature is a specialized subset of the broader account creation capability within the in-house credit card core system for property and casualty insurance personal lines. It enables the full lifecycle management of accounts that have transitioned from good standing to charged-off status due to triggers such as prolonged delinquency, customer death, or bankruptcy. This feature is designed to establish charged-off accounts within the modern core infrastructure, supporting real-time credit issuance, account management, and servicing processes while ensuring seamless integration with credit lines, authorizations, settlements, billing, payments, interest calculations, fees, rewards, disputes, and delinquency workflows."
Replace-Once $old11 $new11

# 5. Feature Overview - second body paragraph text
$old12 = "This feature includes defining and implementing the necessary functionality to create and manage charged-off accounts on the new core platform, replacing legacy third-party systems. It excludes broader account creation activities unrelated to charged-off status and focuses on the specific business rules and data flows associated with charged-off accounts. Key constraints include compliance with financial regulations, accurate interest and fee calculations, and integration with transaction posting, account updates, and dispute management systems. Strategically, this feature supports the organization’s goal of full ownership and control over credit card operations, improving operational efficiency, risk management, and customer servicing capabilities."
$new12 = "This feature includes the creation and management of charged-off accounts but excludes the initial account setup for accounts in good standing or other non-charged-off account types. It leverages data inputs such as transaction records, account updates, and billing/payment details, and outputs updated account statuses and balances. Critical constraints include maintaining compliance with financial regulations, ensuring data integrity during status transitions, and supporting real-time processing performance. Strategically, this feature is essential for achieving full ownership of credit issuance and account servicing, reducing reliance on legacy third-party systems, and enhancing operational control and risk management within the trade credit ecosystem."
Replace-Once $old12 $new12

# 6. Feature Overview - References line wording tweak
#    (By this point the Product Overview "No external sources used." text has
#    already been removed above, so this now uniquely matches the Feature Overview copy.)
$old14 = "No external sources used."
$new14 = "No external sources were used."
Replace-Once $old14 $new14

# Remove the temporary trailing paragraph added above, restoring the original
# paragraph count/structure.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Delete()

